$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Swap A/B widths; widen G-L; add M (offset by +5/6 padding the engine applies to ColumnWidth)
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 46.166666666666664
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666
$ws.Columns.Item(8).ColumnWidth = 14.166666666666666
$ws.Columns.Item(9).ColumnWidth = 15.166666666666666
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 13.166666666666666
$ws.Columns.Item(12).ColumnWidth = 13.166666666666666
$ws.Columns.Item(13).ColumnWidth = 12.166666666666666

# --- New column M: copy header style (bold/border/centered) from L1 before writing values ---
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Organization Name"
$ws.Range("C1").Value = "Organization Link"
$ws.Range("D1").Value = "Logo Link"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Phone Number"
$ws.Range("H1").Value = "Linkedin Link"
$ws.Range("I1").Value = "Instagram Link"
$ws.Range("J1").Value = "Facebook Link"
$ws.Range("K1").Value = "Twitter Link"
$ws.Range("L1").Value = "Youtube Link"
$ws.Range("M1").Value = "Tiktok Link"

# --- Swap columns A and B for data rows (2-21) ---
for ($r = 2; $r -le 21; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $b
    $ws.Cells.Item($r, 2).Value = $a
}
